# Apply local changes after rebase:
# Add a new task row (row 33) to the tasks registry sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

$ws.Range("A$row").Value = "MAN-20260115-001"
$ws.Range("B$row").Value = "MANUAL-20260115"
$ws.Range("C$row").Value = "Praveen"
$ws.Range("D$row").Value = "praveen"
$ws.Range("E$row").Value = "DELETED"
$ws.Range("F$row").Value = "MEDIUM"

# Due Date (date only)
$base = Get-Date -Year 2026 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0
$ws.Range("G$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G$row").Value = $base

$ws.Range("I$row").Value = "Testing"

# Created On / Last Updated timestamps (fractional day precision)
$ws.Range("J$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J$row").Value = $base.AddDays(0.78684002315)

$ws.Range("K$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K$row").Value = $base.AddDays(0.79896162469)
